# DG: updated logic, model, and storage sections
#
# Converts an EMU (English Metric Unit) offset into the point value that,
# once round-tripped through PowerPoint's (single-precision) Left/Top/
# Width/Height properties, lands back on the exact target EMU.
function ToPt($emu) {
    return ($emu + 0.5) / 12700.0
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------------
# 1. Footer "last updated" date fields: 12/5/2018 -> 4/1/2019
#    These live on every slide layout, the slide master, and the notes
#    master (not on the slide itself).
# ---------------------------------------------------------------------
$customLayouts = $p.Designs.Item(1).SlideMaster.CustomLayouts
for ($i = 1; $i -le $customLayouts.Count; $i++) {
    $layout = $customLayouts.Item($i)
    foreach ($shp in $layout.Shapes) {
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = "4/1/2019"
        }
    }
}

foreach ($shp in $p.Designs.Item(1).SlideMaster.Shapes) {
    if ($shp.Name -like "Date Placeholder*") {
        $shp.TextFrame.TextRange.Text = "4/1/2019"
    }
}

foreach ($shp in $p.NotesMaster.Shapes) {
    if ($shp.Name -like "Date Placeholder*") {
        $shp.TextFrame.TextRange.Text = "4/1/2019"
    }
}

# ---------------------------------------------------------------------
# 2. Model section: VersionedAddressBook -> VersionedFinanceTracker
#    (also shrinks the label font from 10.5pt to 10pt)
# ---------------------------------------------------------------------
$sh = $s.Shapes.Item(11)
$sh.TextFrame.TextRange.Text = "VersionedFinanceTracker"
$sh.TextFrame.TextRange.Font.Size = 10

# UniquePersonList -> UniqueRecordList
$s.Shapes.Item(14).TextFrame.TextRange.Text = "UniqueRecordList"

# Person -> Record
$s.Shapes.Item(16).TextFrame.TextRange.Text = "Record"

# ---------------------------------------------------------------------
# 3. Storage / field-list boxes on the right-hand side of the diagram.
# ---------------------------------------------------------------------
# Name -> Amount (widened)
$sh = $s.Shapes.Item(19)
$sh.Left = ToPt(7712396)
$sh.Width = ToPt(812517)
$sh.TextFrame.TextRange.Text = "Amount"

# Phone -> Date (widened)
$sh = $s.Shapes.Item(22)
$sh.Width = ToPt(812518)
$sh.TextFrame.TextRange.Text = "Date"

# connector feeding the (now) "Date" box
$s.Shapes.Item(23).LockAspectRatio = -1

# Email -> Description (widened + shortened)
$sh = $s.Shapes.Item(24)
$sh.Width = ToPt(812519)
$sh.Height = ToPt(279461)
$sh.TextFrame.TextRange.Text = "Description"

# connector feeding the (now) "Description" box
$sh = $s.Shapes.Item(25)
$sh.LockAspectRatio = -1
$sh.Height = ToPt(315034)

# Address -> Category (widened)
$sh = $s.Shapes.Item(26)
$sh.Left = ToPt(7712396)
$sh.Width = ToPt(812519)
$sh.TextFrame.TextRange.Text = "Category"

# connector feeding the (now) "Category" box
$sh = $s.Shapes.Item(27)
$sh.LockAspectRatio = -1
$sh.Width = ToPt(434401)

# ---------------------------------------------------------------------
# 4. ReadOnlyAddressBook interface -> ReadOnlyFinanceTracker (repositioned
#    and resized to fit the longer label)
# ---------------------------------------------------------------------
$sh = $s.Shapes.Item(30)
$sh.Left = ToPt(1119866)
$sh.Top = ToPt(1946417)
$sh.Width = ToPt(1584718)
$sh.Height = ToPt(416329)
$sh.TextFrame.TextRange.Text = "ReadOnlyFinanceTracker"

# ---------------------------------------------------------------------
# 5. Tag -> Name (widened)
# ---------------------------------------------------------------------
$sh = $s.Shapes.Item(39)
$sh.Width = ToPt(812516)
$sh.TextFrame.TextRange.Text = "Name"

# Drop the stray "*" multiplicity marker next to the (now) "Name" box.
$sh = $s.Shapes.Item(41)
$sh.TextFrame.TextRange.Characters(1, 1).Delete()

# AddressBook -> FinanceTracker
$s.Shapes.Item(42).TextFrame.TextRange.Text = "FinanceTracker"
